$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.303.58"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "1.904.72"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.694"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "246.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.71"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.350"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.01"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.24%  "

$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0994"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "2.179.53"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.713"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.60%  "

$ws.Range("D16").Value = "1.906.75"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.85"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.15%  "

$ws.Range("D18").Value = "35.291.41"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "0.0₃0822"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.67%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.44"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.89%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.131"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.62%  "

$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.966"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0573"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.15"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.96%  "

$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0658"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.11"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +2.63%  "

$ws.Range("E43").Value = "  +4.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "90.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "1.347.06"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.50%  "

$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.25%  "

$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.82%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("E51").Value = "  -1.95%  "
